$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5361516.5
$ws.Range("I32").Value = 212
$ws.Range("J32").Value = 7744319
$ws.Range("K32").Value = 212
$ws.Range("L32").Value = 7744319
$ws.Range("M32").Value = 114
$ws.Range("N32").Value = -7744971

$ws.Range("H40").Value = 17483.334
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 17483.334
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 17483.334
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -17833.334

$ws.Range("H43").Value = 88383.92
$ws.Range("I43").Value = 4598.2
$ws.Range("J43").Value = 140750
$ws.Range("K43").Value = 4598.2
$ws.Range("L43").Value = 140750
$ws.Range("M43").Value = -4529.2
$ws.Range("N43").Value = -140888

$ws.Range("H64").Value = 3890.4333
$ws.Range("I64").Value = 3699.8667
$ws.Range("J64").Value = 4081
$ws.Range("K64").Value = 3699.8667
$ws.Range("L64").Value = 4081
$ws.Range("M64").Value = -3451.8667
$ws.Range("N64").Value = -4577

$ws.Range("H67").Value = 3890.4333
$ws.Range("I67").Value = 3699.8667
$ws.Range("J67").Value = 4081
$ws.Range("K67").Value = 3699.8667
$ws.Range("L67").Value = 4081
$ws.Range("M67").Value = -2841.8667
$ws.Range("N67").Value = -5797

$ws.Range("H111").Value = 1645.2
$ws.Range("I111").Value = 1460
$ws.Range("J111").Value = 1768.6666
$ws.Range("K111").Value = 4380
$ws.Range("L111").Value = 5305.9998
$ws.Range("M111").Value = -1313
$ws.Range("N111").Value = -11439.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2291.6667
$ws.Range("I2").Value = 2745.5557
$ws.Range("J2").Value = 930
$ws.Range("K2").Value = 2745.5557
$ws.Range("L2").Value = 930
$ws.Range("M2").Value = -2632.5557
$ws.Range("N2").Value = -1156

$ws.Range("H45").Value = 1528.2632
$ws.Range("I45").Value = 1378.3334
$ws.Range("J45").Value = 1785.2858
$ws.Range("K45").Value = 1378.3334
$ws.Range("L45").Value = 1785.2858
$ws.Range("M45").Value = -1001.3334
$ws.Range("N45").Value = -2539.2858

$ws.Range("H63").Value = 18500
$ws.Range("I63").Value = 21333.334
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 21333.334
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -20647.334
$ws.Range("N63").Value = -11372

$ws.Range("H66").Value = 18500
$ws.Range("I66").Value = 21333.334
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 106666.67
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -103234.67
$ws.Range("N66").Value = -56864

$ws.Range("H74").Value = 1720.4219
$ws.Range("I74").Value = 1034.9073
$ws.Range("J74").Value = 5422.2
$ws.Range("K74").Value = 1034.9073
$ws.Range("L74").Value = 5422.2
$ws.Range("M74").Value = -160.9073000000001
$ws.Range("N74").Value = -7170.2

$ws.Range("H75").Value = 31445
$ws.Range("J75").Value = 31445
$ws.Range("L75").Value = 31445
$ws.Range("N75").Value = -33193

$ws.Range("H77").Value = 1720.4219
$ws.Range("I77").Value = 1034.9073
$ws.Range("J77").Value = 5422.2
$ws.Range("K77").Value = 5174.5365
$ws.Range("L77").Value = 27111
$ws.Range("M77").Value = -806.5365000000002
$ws.Range("N77").Value = -35847

$ws.Range("H78").Value = 31445
$ws.Range("J78").Value = 31445
$ws.Range("L78").Value = 94335
$ws.Range("N78").Value = -103071

$ws.Range("H116").Value = 2291.6667
$ws.Range("I116").Value = 2745.5557
$ws.Range("J116").Value = 930
$ws.Range("K116").Value = 2745.5557
$ws.Range("L116").Value = 930
$ws.Range("M116").Value = -451.5556999999999
$ws.Range("N116").Value = -5518

$ws.Range("H122").Value = 2817.652
$ws.Range("I122").Value = 2961.647
$ws.Range("J122").Value = 2409.6667
$ws.Range("K122").Value = 8884.940999999999
$ws.Range("L122").Value = 7229.000100000001
$ws.Range("M122").Value = -6434.940999999999
$ws.Range("N122").Value = -12129.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2291.6667
$ws.Range("I3").Value = 2745.5557
$ws.Range("J3").Value = 930
$ws.Range("K3").Value = 2745.5557
$ws.Range("L3").Value = 930
$ws.Range("M3").Value = -2631.5557
$ws.Range("N3").Value = -1158

$ws.Range("H20").Value = 2183.4443
$ws.Range("I20").Value = 2211.2307
$ws.Range("J20").Value = 2157.6428
$ws.Range("K20").Value = 2211.2307
$ws.Range("L20").Value = 2157.6428
$ws.Range("M20").Value = -1964.2307
$ws.Range("N20").Value = -2651.6428

$ws.Range("H80").Value = 326.76
$ws.Range("I80").Value = 151
$ws.Range("J80").Value = 409.47058
$ws.Range("K80").Value = 151
$ws.Range("L80").Value = 409.47058
$ws.Range("M80").Value = 847
$ws.Range("N80").Value = -2405.47058

$ws.Range("H83").Value = 326.76
$ws.Range("I83").Value = 151
$ws.Range("J83").Value = 409.47058
$ws.Range("K83").Value = 755
$ws.Range("L83").Value = 2047.3529
$ws.Range("M83").Value = 4237
$ws.Range("N83").Value = -12031.3529

$ws.Range("H99").Value = 6926603
$ws.Range("I99").Value = 2266356.8
$ws.Range("J99").Value = 33334664
$ws.Range("K99").Value = 2266356.8
$ws.Range("L99").Value = 33334664
$ws.Range("M99").Value = -2264858.8
$ws.Range("N99").Value = -33337660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2663.309
$ws.Range("I31").Value = 1885.1621
$ws.Range("J31").Value = 4262.8335
$ws.Range("K31").Value = 1885.1621
$ws.Range("L31").Value = 4262.8335
$ws.Range("M31").Value = -1590.1621
$ws.Range("N31").Value = -4852.8335

$ws.Range("H34").Value = 2663.309
$ws.Range("I34").Value = 1885.1621
$ws.Range("J34").Value = 4262.8335
$ws.Range("K34").Value = 1885.1621
$ws.Range("L34").Value = 4262.8335
$ws.Range("M34").Value = -1683.1621
$ws.Range("N34").Value = -4666.8335

$ws.Range("H86").Value = 5715.1177
$ws.Range("I86").Value = 3911.9
$ws.Range("J86").Value = 8291.143
$ws.Range("K86").Value = 3911.9
$ws.Range("L86").Value = 8291.143
$ws.Range("M86").Value = -2788.9
$ws.Range("N86").Value = -10537.143

$ws.Range("H89").Value = 5715.1177
$ws.Range("I89").Value = 3911.9
$ws.Range("J89").Value = 8291.143
$ws.Range("K89").Value = 19559.5
$ws.Range("L89").Value = 41455.715
$ws.Range("M89").Value = -13943.5
$ws.Range("N89").Value = -52687.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 585.3182
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 637.35
$ws.Range("K23").Value = 195
$ws.Range("L23").Value = 1912.05
$ws.Range("M23").Value = 40
$ws.Range("N23").Value = -2382.05

$ws.Range("H122").Value = 587.875
$ws.Range("I122").Value = 382.27777
$ws.Range("J122").Value = 852.2143
$ws.Range("K122").Value = 3440.49993
$ws.Range("L122").Value = 7669.928699999999
$ws.Range("M122").Value = -990.4999299999999
$ws.Range("N122").Value = -12569.9287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7294.028
$ws.Range("I70").Value = 3826.2
$ws.Range("J70").Value = 24633.166
$ws.Range("K70").Value = 3826.2
$ws.Range("L70").Value = 24633.166
$ws.Range("M70").Value = -3556.2
$ws.Range("N70").Value = -25173.166

$ws.Range("H73").Value = 7294.028
$ws.Range("I73").Value = 3826.2
$ws.Range("J73").Value = 24633.166
$ws.Range("K73").Value = 3826.2
$ws.Range("L73").Value = 24633.166
$ws.Range("M73").Value = -2890.2
$ws.Range("N73").Value = -26505.166

$ws.Range("H132").Value = 4646.0625
$ws.Range("I132").Value = 5267.5293
$ws.Range("J132").Value = 3941.7334
$ws.Range("K132").Value = 15802.5879
$ws.Range("L132").Value = 11825.2002
$ws.Range("M132").Value = -13272.5879
$ws.Range("N132").Value = -16885.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3570
$ws.Range("I16").Value = 910
$ws.Range("J16").Value = 4900
$ws.Range("K16").Value = 910
$ws.Range("L16").Value = 4900
$ws.Range("M16").Value = -740
$ws.Range("N16").Value = -5240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2309.756
$ws.Range("I132").Value = 1135.8846
$ws.Range("J132").Value = 4344.467
$ws.Range("K132").Value = 3407.6538
$ws.Range("L132").Value = 13033.401
$ws.Range("M132").Value = -877.6538
$ws.Range("N132").Value = -18093.401
